# Regenerate save_data: replace column G ("K") values (previously derived
# from "Strike#") for rows 2-29 with the newly computed K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 2
    4  = 2
    5  = 3
    6  = 8
    7  = 8
    8  = 2
    9  = 6
    10 = 3
    11 = 4
    12 = 4
    13 = 3
    14 = 2
    15 = 6
    16 = 6
    17 = 8
    18 = 2
    19 = 6
    20 = 2
    21 = 5
    22 = 1
    23 = 7
    24 = 4
    25 = 2
    26 = 2
    27 = 3
    28 = 2
    29 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
